$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24 (shifts existing rows 24.. down by one)
$ws.Rows.Item(24).Insert()

# Populate the newly inserted row 24 with the new record
$ws.Cells.Item(24, 1).Value = 7
$ws.Cells.Item(24, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(24, 3).Value = "Ñuble"
$ws.Cells.Item(24, 4).Value = 44949
$ws.Cells.Item(24, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(24, 5).Value = 16
$ws.Cells.Item(24, 6).Value = "Fruta"
$ws.Cells.Item(24, 7).Value = 100108
$ws.Cells.Item(24, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(24, 9).Value = 100108002
$ws.Cells.Item(24, 10).Value = "Mango"
$ws.Cells.Item(24, 11).Value = "Sin especificar"
$ws.Cells.Item(24, 12).Value = "Primera"
$ws.Cells.Item(24, 13).Value = 50
$ws.Cells.Item(24, 14).Value = 7000
$ws.Cells.Item(24, 15).Value = 7000
$ws.Cells.Item(24, 16).Value = 7000
$ws.Cells.Item(24, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(24, 18).Value = "Perú"
$ws.Cells.Item(24, 19).Value = 1750
$ws.Cells.Item(24, 20).Value = 4
